# Apply the "base" macro2 workbook change:
#  - add new base command outputToCloud(resource) to the 'base' list (column E)
#  - add new named-range category 'text' (a single new column Y) holding spellCheck(var,profile,text)
#  - add 'text' to the 'target' list (column A)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# 1) Insert 'text' into the sorted "target" list (column A), which is used by
#    the 'target' defined name ($A$2:$A$30 -> $A$2:$A$31).
#    'text' sorts alphabetically right after 'step' and before 'web', i.e.
#    at row 25 (pushing the previous A25:A30 down to A26:A31).
# ---------------------------------------------------------------------------
for ($r = 30; $r -ge 25; $r--) {
    $ws.Cells.Item($r + 1, 1).Value2 = $ws.Cells.Item($r, 1).Value2
}
$ws.Cells.Item(25, 1).Value2 = "text"

# ---------------------------------------------------------------------------
# 2) Insert 'outputToCloud(resource)' into the sorted "base" list (column E),
#    used by the 'base' defined name ($E$2:$E$38 -> $E$2:$E$39).
#    It sorts alphabetically right after 'macro(file,sheet,name)' and before
#    'prependText(var,prependWith)', i.e. at row 22 (pushing E22:E38 down to
#    E23:E39).
# ---------------------------------------------------------------------------
for ($r = 38; $r -ge 22; $r--) {
    $ws.Cells.Item($r + 1, 5).Value2 = $ws.Cells.Item($r, 5).Value2
}
$ws.Cells.Item(22, 5).Value2 = "outputToCloud(resource)"

# ---------------------------------------------------------------------------
# 3) Insert a brand-new column at Y to host the new 'text' named-range
#    category. This shifts the existing Y:AD columns (web, webalert,
#    webcookie, ws, ws.async, xml) one column to the right (Z:AE).
# ---------------------------------------------------------------------------
$ws.Columns("Y:Y").Insert()

$ws.Cells.Item(1, 25).Value2 = "text"
$ws.Cells.Item(2, 25).Value2 = "spellCheck(var,profile,text)"

# ---------------------------------------------------------------------------
# 4) Update defined names so that ranges reflect the new layout.
# ---------------------------------------------------------------------------
$wb.Names.Item("base").RefersTo = "='#system'!`$E`$2:`$E`$39"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$31"
$wb.Names.Item("web").RefersTo = "='#system'!`$Z`$2:`$Z`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AC`$2:`$AC`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AD`$2:`$AD`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AE`$2:`$AE`$27"
$wb.Names.Add("text", "='#system'!`$Y`$2:`$Y`$2")
